# Apply the author's edit: the daily-report sheet had its "cost / sold /
# margin" columns (C:D:E, rows 1-108) reformatted from a 2-decimal custom
# number format to a 3-decimal one, and left selected when the file was
# saved.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("C1:E108")
$dataRange.NumberFormat = "[$-1010419]#,##0.000;\-#,##0.000"

# Leave the just-edited range selected, matching the saved view state.
$dataRange.Select()
